# Generate Report for Handoff
# Adds a new localization-status row (for file
# f111199a-0372-4e00-97d4-5cd340f5d427.md) to the Overview, zh-cn and
# de-de sheets, mirroring the existing 8a0554f4-... row.

$wb = $excel.ActiveWorkbook

$newFile        = "f111199a-0372-4e00-97d4-5cd340f5d427.md"
$newPath        = "e2e\f111199a-0372-4e00-97d4-5cd340f5d427.md"
$hyperlinkUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d79979c5d3418168546c54e4900311eae1f48318/e2e/f111199a-0372-4e00-97d4-5cd340f5d427.md"
$dateFormat     = "yyyy-mm-dd HH:mm:ss"
$hyperlinkColor = 15570276   # BGR for RGB FF6495ED, matches existing HyperLink style

# ---------------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item("Overview")
$tblOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $newPath
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $newPath) | Out-Null
$wsOverview.Range("B3").Font.Underline = 1
$wsOverview.Range("B3").Font.Color = $hyperlinkColor
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").NumberFormat = $dateFormat
$wsOverview.Range("G3").Value = "2016-08-16 06:35:17"

# ---------------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$tblZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$tblZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkUrl, "", "", $newFile) | Out-Null
$wsZhCn.Range("A3").Font.Underline = 1
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "f111199a-0372-4e00-97d4-5cd340f5d427.8f13ce76e796dba86417cca1c1795cdcc6d0dbf9.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("H3").Value = "2016-08-16 06:35:12"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = ""

# ---------------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$tblDeDe = $wsDeDe.ListObjects.Item("de-de")
$tblDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkUrl, "", "", $newFile) | Out-Null
$wsDeDe.Range("A3").Font.Underline = 1
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "f111199a-0372-4e00-97d4-5cd340f5d427.8f13ce76e796dba86417cca1c1795cdcc6d0dbf9.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("H3").Value = "2016-08-16 06:35:17"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = ""
